$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (i.e. right
#    before "总计"), copying its formatting from "2021-Q4" so fonts,
#    borders and alignment match the rest of the workbook.
# -------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

$q4Sheet.Range("A1:H13").Copy()
$newSheet.Range("A1:H13").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns that must stay TEXT (fund code + the numeric-looking text
# columns that were authored as text in the source data).
$newSheet.Range("B2:B13").NumberFormat = "@"
$newSheet.Range("D2:G13").NumberFormat = "@"

# Row index column (A) - plain numbers 0..11
$indexValues = @(0,1,2,3,4,5,6,7,8,9,10,11)
for ($i = 0; $i -lt $indexValues.Count; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $indexValues[$i]
}

# Fund holdings data: code, name, fund size, total stock position,
# position ratio, held market value (100M yuan), position rank
$rows = @(
    @("502023", "鹏华国证钢铁行业指数（LOF）", "15.55", "94.76", "11.24", "1.7478", 2),
    @("515210", "国泰中证钢铁ETF", "16.24", "99.25", "8.53", "1.3853", 2),
    @("516150", "嘉实中证稀土产业ETF", "25.17", "99.75", "5.12", "1.2887", 5),
    @("002666", "前海开源沪港深创新成长灵活配置混合A", "11.96", "81.64", "7.84", "0.9377", 3),
    @("516780", "华泰柏瑞中证稀土产业ETF", "11.06", "98.70", "5.12", "0.5663", 5),
    @("168203", "中融国证钢铁行业指数", "4.30", "92.58", "10.93", "0.4700", 2),
    @("002667", "前海开源沪港深创新成长灵活配置混合C", "3.25", "81.64", "7.84", "0.2548", 3),
    @("159715", "易方达中证稀土产业ETF", "3.42", "99.06", "5.08", "0.1737", 5),
    @("159713", "富国中证稀土产业交易型开放式指数证券投资基金", "3.26", "99.26", "5.11", "0.1666", 5),
    @("013802", "财通资管中证钢铁指数A", "0.11", "90.83", "7.74", "0.0085", 2),
    @("159944", "广发中证全指原材料ETF", "0.20", "98.36", "1.14", "0.0023", 9),
    @("013803", "财通资管中证钢铁指数C", "0.02", "90.83", "7.74", "0.0015", 2)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# -------------------------------------------------------------------------
# 2. Add a "2022-Q1" summary row at the top of the "总计" sheet's data
#    (right after the header row), pushing all the other rows down and
#    renumbering the index column.
# -------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

$totalWs.Rows.Item(2).Insert()
$totalWs.Range("A2:D2").ClearFormats()

# Re-apply the index-column style (bold, bordered, centered) that the
# rest of column A uses.
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 12
$totalWs.Range("D2").Value = 7

# Renumber the remaining index values (they shifted down by one row).
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
$totalWs.Range("A7").Value = 5

Write-Host "2022-Q1 sheet added and 总计 sheet updated."
